$wb = $excel.ActiveWorkbook

# --- Summary ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = "2026-01-18 13:58 UTC"
$ws.Range("B4").Value = 53
$ws.Range("B5").Value = 53
$ws.Range("B6").Value = 53

# --- Reversal Setups ---
$ws = $wb.Worksheets.Item("Reversal Setups")
$ws.Range("B2").Value = "EGLDUSDT"
$ws.Range("C2").Value = "MultiversX"
$ws.Range("D2").Value = "`$6.55"
$ws.Range("E2").Value = "`$191.42M"
$ws.Range("F2").Value = "`$2.09M"
$ws.Range("G2").Value = 76.95
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 59.73
$ws.Range("B3").Value = "MEUSDT"
$ws.Range("C3").Value = "Magic Eden"
$ws.Range("D3").Value = "`$0.29"
$ws.Range("E3").Value = "`$124.94M"
$ws.Range("F3").Value = "`$4.03M"
$ws.Range("G3").Value = 75
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 100
$ws.Range("B4").Value = "MANAUSDT"
$ws.Range("C4").Value = "Decentraland"
$ws.Range("D4").Value = "`$0.16"
$ws.Range("E4").Value = "`$312.37M"
$ws.Range("F4").Value = "`$3.41M"
$ws.Range("G4").Value = 75
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("B5").Value = "GALAUSDT"
$ws.Range("C5").Value = "Gala"
$ws.Range("D5").Value = "`$0.01"
$ws.Range("E5").Value = "`$355.54M"
$ws.Range("F5").Value = "`$2.74M"
$ws.Range("G5").Value = 72.5
$ws.Range("I5").Value = 80
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 0
$ws.Range("B6").Value = "SANDUSDT"
$ws.Range("C6").Value = "The Sandbox"
$ws.Range("D6").Value = "`$0.15"
$ws.Range("E6").Value = "`$389.86M"
$ws.Range("F6").Value = "`$4.68M"
$ws.Range("G6").Value = 70.44
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 77.20999999999999
$ws.Range("B7").Value = "CAKEUSDT"
$ws.Range("C7").Value = "PancakeSwap"
$ws.Range("D7").Value = "`$2.08"
$ws.Range("E7").Value = "`$714.90M"
$ws.Range("F7").Value = "`$1.62M"
$ws.Range("G7").Value = 67.5
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 50
$ws.Range("B8").Value = "BERAUSDT"
$ws.Range("C8").Value = "Berachain"
$ws.Range("D8").Value = "`$0.87"
$ws.Range("E8").Value = "`$125.92M"
$ws.Range("F8").Value = "`$6.28M"
$ws.Range("G8").Value = 65.76000000000001
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 53.79
$ws.Range("B9").Value = "FILUSDT"
$ws.Range("C9").Value = "Filecoin"
$ws.Range("D9").Value = "`$1.53"
$ws.Range("E9").Value = "`$1.13B"
$ws.Range("F9").Value = "`$3.90M"
$ws.Range("B10").Value = "IMXUSDT"
$ws.Range("C10").Value = "Immutable"
$ws.Range("D10").Value = "`$0.29"
$ws.Range("E10").Value = "`$572.67M"
$ws.Range("F10").Value = "`$1.26M"
$ws.Range("G10").Value = 65
$ws.Range("J10").Value = 60
$ws.Range("B11").Value = "BONKUSDT"
$ws.Range("C11").Value = "Bonk"
$ws.Range("D11").Value = "`$0.00"
$ws.Range("E11").Value = "`$910.13M"
$ws.Range("F11").Value = "`$1.03M"
$ws.Range("G11").Value = 65
$ws.Range("J11").Value = 60

# --- Breakout Setups ---
$ws = $wb.Worksheets.Item("Breakout Setups")
$ws.Range("B2").Value = "SANDUSDT"
$ws.Range("C2").Value = "The Sandbox"
$ws.Range("D2").Value = "`$0.15"
$ws.Range("E2").Value = "`$389.86M"
$ws.Range("F2").Value = "`$4.68M"
$ws.Range("G2").Value = 65
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("B3").Value = "MEUSDT"
$ws.Range("C3").Value = "Magic Eden"
$ws.Range("D3").Value = "`$0.29"
$ws.Range("E3").Value = "`$124.94M"
$ws.Range("F3").Value = "`$4.03M"
$ws.Range("G3").Value = 65
$ws.Range("I3").Value = 100
$ws.Range("B4").Value = "MANAUSDT"
$ws.Range("C4").Value = "Decentraland"
$ws.Range("D4").Value = "`$0.16"
$ws.Range("E4").Value = "`$312.37M"
$ws.Range("F4").Value = "`$3.41M"
$ws.Range("G4").Value = 59.69
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 64.59999999999999
$ws.Range("D5").Value = "`$0.87"
$ws.Range("E5").Value = "`$125.92M"
$ws.Range("F5").Value = "`$6.28M"
$ws.Range("G5").Value = 59.2
$ws.Range("I5").Value = 80.68000000000001
$ws.Range("B6").Value = "METUSDT"
$ws.Range("C6").Value = "Meteora"
$ws.Range("D6").Value = "`$0.29"
$ws.Range("E6").Value = "`$145.40M"
$ws.Range("F6").Value = "`$3.05M"
$ws.Range("G6").Value = 57.23
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 61.5
$ws.Range("B7").Value = "AXSUSDT"
$ws.Range("C7").Value = "Axie Infinity"
$ws.Range("D7").Value = "`$1.87"
$ws.Range("E7").Value = "`$318.20M"
$ws.Range("F7").Value = "`$26.86M"
$ws.Range("G7").Value = 55.79
$ws.Range("I7").Value = 69.29000000000001
$ws.Range("B8").Value = "IPUSDT"
$ws.Range("C8").Value = "Story"
$ws.Range("D8").Value = "`$2.67"
$ws.Range("E8").Value = "`$936.31M"
$ws.Range("F8").Value = "`$5.15M"
$ws.Range("G8").Value = 39.95
$ws.Range("I8").Value = 24.17
$ws.Range("K8").Value = 84.68000000000001
$ws.Range("B9").Value = "EGLDUSDT"
$ws.Range("C9").Value = "MultiversX"
$ws.Range("D9").Value = "`$6.55"
$ws.Range("E9").Value = "`$191.42M"
$ws.Range("F9").Value = "`$2.09M"
$ws.Range("G9").Value = 38.68
$ws.Range("I9").Value = 89.59999999999999
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 25.37
$ws.Range("B10").Value = "DASHUSDT"
$ws.Range("C10").Value = "Dash"
$ws.Range("D10").Value = "`$84.76"
$ws.Range("E10").Value = "`$1.06B"
$ws.Range("F10").Value = "`$28.01M"
$ws.Range("B11").Value = "CHZUSDT"
$ws.Range("C11").Value = "Chiliz"
$ws.Range("D11").Value = "`$0.06"
$ws.Range("E11").Value = "`$613.44M"
$ws.Range("F11").Value = "`$21.44M"

# --- Pullback Setups ---
$ws = $wb.Worksheets.Item("Pullback Setups")
$ws.Range("B2").Value = "MEUSDT"
$ws.Range("C2").Value = "Magic Eden"
$ws.Range("D2").Value = "`$0.29"
$ws.Range("E2").Value = "`$124.94M"
$ws.Range("F2").Value = "`$4.03M"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 90
$ws.Range("K2").Value = 100
$ws.Range("B3").Value = "METUSDT"
$ws.Range("C3").Value = "Meteora"
$ws.Range("D3").Value = "`$0.29"
$ws.Range("E3").Value = "`$145.40M"
$ws.Range("F3").Value = "`$3.05M"
$ws.Range("G3").Value = 70.5
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 100
$ws.Range("B4").Value = "AXSUSDT"
$ws.Range("C4").Value = "Axie Infinity"
$ws.Range("D4").Value = "`$1.87"
$ws.Range("E4").Value = "`$318.20M"
$ws.Range("F4").Value = "`$26.86M"
$ws.Range("G4").Value = 63.5
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 80
$ws.Range("B5").Value = "SANDUSDT"
$ws.Range("C5").Value = "The Sandbox"
$ws.Range("D5").Value = "`$0.15"
$ws.Range("E5").Value = "`$389.86M"
$ws.Range("F5").Value = "`$4.68M"
$ws.Range("G5").Value = 61.5
$ws.Range("K5").Value = 100
$ws.Range("B6").Value = "BERAUSDT"
$ws.Range("C6").Value = "Berachain"
$ws.Range("D6").Value = "`$0.87"
$ws.Range("E6").Value = "`$125.92M"
$ws.Range("F6").Value = "`$6.28M"
$ws.Range("G6").Value = 57.5
$ws.Range("H6").Value = 80
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 80
$ws.Range("B7").Value = "MANAUSDT"
$ws.Range("C7").Value = "Decentraland"
$ws.Range("D7").Value = "`$0.16"
$ws.Range("E7").Value = "`$312.37M"
$ws.Range("F7").Value = "`$3.41M"
$ws.Range("G7").Value = 56.5
$ws.Range("H7").Value = 80
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 100
$ws.Range("B8").Value = "DASHUSDT"
$ws.Range("C8").Value = "Dash"
$ws.Range("D8").Value = "`$84.76"
$ws.Range("E8").Value = "`$1.06B"
$ws.Range("F8").Value = "`$28.01M"
$ws.Range("G8").Value = 55
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = 80
$ws.Range("D9").Value = "`$26.43"
$ws.Range("E9").Value = "`$521.38M"
$ws.Range("F9").Value = "`$5.03M"
$ws.Range("G9").Value = 55
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 0
$ws.Range("B10").Value = "ACHUSDT"
$ws.Range("C10").Value = "Alchemy Pay"
$ws.Range("D10").Value = "`$0.01"
$ws.Range("E10").Value = "`$118.97M"
$ws.Range("F10").Value = "`$1.20M"
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = 100
$ws.Range("J10").Value = 80
$ws.Range("B11").Value = "ZENUSDT"
$ws.Range("C11").Value = "Horizen"
$ws.Range("D11").Value = "`$12.47"
$ws.Range("E11").Value = "`$219.52M"
$ws.Range("F11").Value = "`$5.80M"
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 60
